$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1060
$ws1.Range("F8").Value = 206
$ws1.Range("F10").Value = 2
$ws1.Range("F11").Value = 7
$ws1.Range("F12").Value = 496
$ws1.Range("F15").Value = 12455
$ws1.Range("F16").Value = 129
$ws1.Range("F17").Value = 5492

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 121

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 121
$ws4.Range("F7").Value = 1060
$ws4.Range("F10").Value = 206
$ws4.Range("F12").Value = 2
$ws4.Range("F13").Value = 7
$ws4.Range("F14").Value = 496
$ws4.Range("F17").Value = 12455
$ws4.Range("F19").Value = 129
$ws4.Range("F20").Value = 5492
